$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-03-18 Tuesday"; new="2025-03-19 Wednesday"},
    @{old="582×5="; new="177×9="},
    @{old="705×5="; new="746×5="},
    @{old="440×6="; new="388×3="},
    @{old="141×8="; new="928×2="},
    @{old="384×4="; new="547×3="},
    @{old="965×3="; new="922×2="},
    @{old="732×3="; new="573×8="},
    @{old="856×3="; new="453×4="},
    @{old="694×7="; new="963×8="},
    @{old="209×3="; new="373×2="},
    @{old="612×4="; new="219×4="},
    @{old="875×5="; new="413×6="},
    @{old="625×5="; new="171×2="},
    @{old="261×5="; new="629×2="},
    @{old="246×2="; new="178×8="},
    @{old="353×2="; new="618×4="},
    @{old="313×9="; new="831×9="},
    @{old="574×3="; new="762×4="},
    @{old="452×5="; new="956×9="},
    @{old="297×7="; new="734×9="},
    @{old="796×4="; new="247×2="},
    @{old="927×5="; new="314×2="},
    @{old="364×8="; new="511×2="},
    @{old="991×9="; new="651×7="},
    @{old="260×8="; new="575×8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
